# Insert a new weekly record at row 280 for "Vega Modelo de Temuco" / "Betarraga".
# This shifts all existing rows from 280 down to 281 (old row 280 -> 281, ...,
# old row 360 -> 361), growing the used range from A1:R360 to A1:R361, then
# fills the newly-inserted row 280 with the new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 280:360 down to 281:361, leaving a blank row 280 behind.
$ws.Rows.Item(280).Insert()

# Populate the new row 280 with the new record.
$ws.Cells.Item(280, 1).Value2  = 10
$ws.Cells.Item(280, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(280, 3).Value2  = "La Araucanía"
$ws.Cells.Item(280, 4).Value2  = 44627
$ws.Cells.Item(280, 5).Value2  = 9
$ws.Cells.Item(280, 6).Value2  = 100114014
$ws.Cells.Item(280, 7).Value2  = "Betarraga"
$ws.Cells.Item(280, 8).Value2  = "Sin especificar"
$ws.Cells.Item(280, 9).Value2  = "Primera"
$ws.Cells.Item(280, 10).Value2 = 110
$ws.Cells.Item(280, 11).Value2 = 10000
$ws.Cells.Item(280, 12).Value2 = 10000
$ws.Cells.Item(280, 13).Value2 = 10000
$ws.Cells.Item(280, 14).Value2 = "`$/docena de paquetes"
$ws.Cells.Item(280, 15).Value2 = "Región del Maule"
$ws.Cells.Item(280, 16).Value2 = 833
$ws.Cells.Item(280, 17).Value2 = 12
$ws.Cells.Item(280, 18).Value2 = "Hortaliza"
